$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 was the text label "Parameter_Name" - change it to the literal number 3
$ws.Range("A1").Value = 3

# Apply Text number format ("@") to the column-B data-entry cells so that
# values like serial numbers / grating numbers are stored as text instead
# of General (this also flips the column-C default/A18 style bookkeeping
# automatically once the old unused style is no longer referenced).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
